$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1099-CAP")

# Update the header for column V
$ws.Range("V1").Value = "Is Corrected Form of 1099"

# Update the "Is Corrected" values for each data row to reflect corrected status
$ws.Range("V2").Value = "Yes"
$ws.Range("V3").Value = "Yes"
$ws.Range("V4").Value = "Yes"
$ws.Range("V5").Value = "Yes"
$ws.Range("V6").Value = "Yes"
$ws.Range("V7").Value = "Yes"
$ws.Range("V8").Value = "No"
$ws.Range("V9").Value = "No"
$ws.Range("V10").Value = "No"
$ws.Range("V11").Value = "No"
$ws.Range("V12").Value = "Yes"
$ws.Range("V13").Value = "Yes"

# Update the active selection to V13
$ws.Range("V13").Select()
